# Applies the "adding experiment and animal to excel" edit:
#  - adds "animal name"/"experiment" (D/E) columns with data to the
#    Example3, Example2 and Example1 sheets
#  - updates the view/selection state of every sheet, and makes
#    "Example1" the active sheet/tab

$wb = $excel.ActiveWorkbook

$workingSheet = $wb.Worksheets.Item(1)   # "working sheet"
$example3     = $wb.Worksheets.Item(2)   # "Example3"
$example2     = $wb.Worksheets.Item(3)   # "Example2"
$example1     = $wb.Worksheets.Item(4)   # "Example1"

# --- "working sheet" -------------------------------------------------
# Only the selection/view changes here; the animal/experiment headers
# already exist in D1/E1 on this sheet.
$workingSheet.Range("D1:E1").Select()

# --- Example3 ----------------------------------------------------------
$example3.Range("D1").Value = "animal name"
$example3.Range("E1").Value = "experiment"

$example3.Range("D2").Value = "M26"
$example3.Range("D3").Value = "M26"
$example3.Range("D4").Value = "M26"

$example3.Range("E3").Value = "9_24_17"
$example3.Range("E2").Value = "9_28_17"
$example3.Range("E4").Value = "10_1_17"

$example3.Range("D2:E4").Select()

# --- Example2 ------------------------------------------------------------
$example2.Range("D1").Value = "animal name"
$example2.Range("E1").Value = "experiment"

$example2.Range("D2").Value = "M26"
$example2.Range("D3").Value = "M26"
$example2.Range("D4").Value = "M26"

$example2.Range("E3").Value = "9_24_17"
$example2.Range("E2").Value = "9_28_17"
$example2.Range("E4").Value = "10_1_17"

$example2.Range("D2:E4").Select()

# --- Example1 (edited/selected last so it ends up the active tab) -------
$example1.Range("D1").Value = "animal name"
$example1.Range("E1").Value = "experiment"

$example1.Range("D2").Value = "M26"
$example1.Range("D3").Value = "M26"
$example1.Range("D4").Value = "M26"

$example1.Range("E3").Value = "9_24_17"
$example1.Range("E2").Value = "9_28_17"
$example1.Range("E4").Value = "10_1_17"

$example1.Range("D15").Select()
